# Timesheet changes by Ruchika
#
# For the week rows 28-31 ("FebruaryMarch 2013" sheet), columns AP:AU were
# blank (weekend columns that should carry the same "OFF" marker already
# present in column AO). Fill AP:AU with the same "OFF" value/format that
# column AO already uses, for each of the four rows.
#
# Also update the active sheet's selection to match the edited area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FebruaryMarch 2013")
$ws.Activate()

$rows = 28..31
$cols = @("AP", "AQ", "AR", "AS", "AT", "AU")

foreach ($r in $rows) {
    # AO<r> already holds the "OFF" shared-string value with the correct
    # style (s="20"); copy both value and formatting across to AP:AU.
    $src = $ws.Range("AO" + $r)
    foreach ($col in $cols) {
        $dst = $ws.Range($col + $r)
        $src.Copy($dst)
    }
}

# Reflect the new selection / active cell left behind by the edit.
$ws.Range("AU28:AU31").Select()
